# Rename the placeholder "Maersk 172".."Maersk 189" vessels that previously
# held individual vessel names (class Aegir / Equinox) and switch their
# PROJECT code from "MM868" to "M868".
#
# Column layout on Sheet1: A=NAME, B=IMO, C=CLASS, D=PROJECT
# Rows 398-415 are the affected rows (CLASS stays Aegir/Equinox, IMO stays
# the same; only NAME and PROJECT change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 398
$lastRow  = 415
$startNum = 172

# Write all the new NAME values first, then all the new PROJECT values, so
# that newly-introduced shared strings land in the same order as authored
# (names before the project code) when the workbook is saved.
$n = $startNum
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "Maersk $n"
    $n = $n + 1
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "M868"
}

$ws.Range("F411").Select()
